$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2..H2): category/weight change to "Техника" / "до 2 кг";
# new "Размер/Комплектация" value "24242" (must stay text, not become a number);
# Артикул -> "fsf"; Сумма -> 4671.
$ws.Range("C2").Value = "Техника"
$ws.Range("D2").Value = "до 2 кг"

# "24242" looks numeric, so Excel would silently convert a plain .Value
# assignment into a number. Force it through a text-formatted helper cell,
# then paste only the *value* into E2 so E2 itself keeps the sheet's
# original (default) style while still carrying a text cell value.
$ws.Range("ZZ100").NumberFormat = "@"
$ws.Range("ZZ100").Value = "24242"
$ws.Range("ZZ100").Copy()
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("ZZ100").Clear()

$ws.Range("F2").Value = "fsf"
$ws.Range("G2").Value = 4671

# Row 3 (A3..H3): category change to "Аксессуары"; weight changes to "до 1 кг";
# Размер/Комплектация stays "-"; Артикул -> "fwgw"; Сумма -> 3480.
$ws.Range("C3").Value = "Аксессуары"
$ws.Range("D3").Value = "до 1 кг"
$ws.Range("F3").Value = "fwgw"
$ws.Range("G3").Value = 3480

$wb.Save()
